$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so moved rows don't leave stale cells behind.
$ws.UsedRange.Clear()

# The sheet was renamed from "Data" to "Summary".
$ws.Name = "Summary"

# Row 1 - country name
$ws.Range("A1").Value = "Bosnia and Herzegovina"
$ws.Range("A1").Font.Size = 18

# Row 3 - section title
$ws.Range("A3").Value = "MSME Participation on the Economy"
$ws.Range("A3").Font.Bold = $true

# Row 8 - new source-type line (bold + underlined)
$ws.Range("A8").Value = "Source Type: Ministry of Finance/Central Bank (Most Widely Used)"
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Underline = $true

# Row 10 - column headers
$ws.Range("B10").Value = "Micro"
$ws.Range("B10").Font.Bold = $true
$ws.Range("C10").Value = "SMEs"
$ws.Range("C10").Font.Bold = $true
$ws.Range("D10").Value = "MSMEs"
$ws.Range("D10").Font.Bold = $true

# Row 11 - Enterprises (absolute #)
$ws.Range("A11").Value = "Enterprises (absolute #)"
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "151107"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "10188"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "161295"

# Row 12 - Enterprises density (per 1000 people)
$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "39.1"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "2.6"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.8"

# Row 13 - new Employment (% of total) row
$ws.Range("A13").Value = "Employment (% of total)"
$ws.Range("A13").Font.Bold = $true
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.8"

# Row 14 - Enterprises (% of total)
$ws.Range("A14").Value = "Enterprises (% of total)"
$ws.Range("A14").Font.Bold = $true
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "93.3"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "6.3"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99.6"

# Row 15 - source citation (italic)
$ws.Range("A15").Value = "Source: MVTEO, 2008"
$ws.Range("A15").Font.Italic = $true

# Row 20 - sector distribution details heading
$ws.Range("A20").Value = "Sector Distribution Details"
$ws.Range("A20").Font.Bold = $true

# Row 23 - new MVTEO heading
$ws.Range("A23").Value = "MVTEO"
$ws.Range("A23").Font.Bold = $true

# Row 24 - new MVTEO citation (italic)
$ws.Range("A24").Value = "Ministry of Foreign Trade and Economic Relations of Bosnia and Herzegovina (MVTEO), ""Small and Medium-Sized Enterprise Development Strategy in Bosnia and Herzegovina 2009 - 2011"", 2009. Available at http://www.mvteo.gov.ba/vijesti/posljednje_vijesti/default.aspx?id=1204&langTag=bs-BA"
$ws.Range("A24").Font.Italic = $true
